$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 628
$ws1.Range("F3").Value = 207
$ws1.Range("F4").Value = 628
$ws1.Range("F5").Value = 552
$ws1.Range("F6").Value = 305
$ws1.Range("F7").Value = 2755
$ws1.Range("F8").Value = 468
$ws1.Range("F9").Value = 7699
$ws1.Range("F12").Value = 37
$ws1.Range("F13").Value = 316
$ws1.Range("F14").Value = 45

# Sheet "全部类型" (All Types) - same underlying events, update column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 628
$ws4.Range("F3").Value = 207
$ws4.Range("F4").Value = 628
$ws4.Range("F5").Value = 552
$ws4.Range("F6").Value = 305
$ws4.Range("F9").Value = 2755
$ws4.Range("F10").Value = 468
$ws4.Range("F11").Value = 7699
$ws4.Range("F14").Value = 37
$ws4.Range("F17").Value = 316
$ws4.Range("F18").Value = 45
